$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H32").Value = 8440.4
$ws.Range("J32").Value = 9925.5
$ws.Range("L32").Value = 9925.5
$ws.Range("N32").Value = -10577.5
$ws.Range("H119").Value = 900
$ws.Range("J119").Value = 900
$ws.Range("L119").Value = 2700
$ws.Range("N119").Value = -12376

$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H74").Value = 4978537
$ws.Range("I74").Value = 6670129
$ws.Range("J74").Value = 3267.353
$ws.Range("K74").Value = 6670129
$ws.Range("L74").Value = 3267.353
$ws.Range("M74").Value = -6669255
$ws.Range("N74").Value = -5015.353
$ws.Range("H77").Value = 4978537
$ws.Range("I77").Value = 6670129
$ws.Range("J77").Value = 3267.353
$ws.Range("K77").Value = 33350645
$ws.Range("L77").Value = 16336.765
$ws.Range("M77").Value = -33346277
$ws.Range("N77").Value = -25072.765
$ws.Range("H102").Value = 1541.5714
$ws.Range("I102").Value = 1429.3846
$ws.Range("K102").Value = 1429.3846
$ws.Range("M102").Value = 192.6153999999999
$ws.Range("H122").Value = 2744.9312
$ws.Range("I122").Value = 1716.421
$ws.Range("J122").Value = 4699.1
$ws.Range("K122").Value = 5149.263
$ws.Range("L122").Value = 14097.3
$ws.Range("M122").Value = -2699.263
$ws.Range("N122").Value = -18997.3
$ws.Range("H132").Value = 3995.9546
$ws.Range("I132").Value = 2679.3333
$ws.Range("K132").Value = 8037.999899999999
$ws.Range("M132").Value = -5507.999899999999

$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H99").Value = 1750.0714
$ws.Range("I99").Value = 1625.1666
$ws.Range("J99").Value = 2499.5
$ws.Range("K99").Value = 1625.1666
$ws.Range("L99").Value = 2499.5
$ws.Range("M99").Value = -127.1666
$ws.Range("N99").Value = -5495.5
$ws.Range("H107").Value = 2459.8
$ws.Range("I107").Value = 599.6667
$ws.Range("K107").Value = 599.6667
$ws.Range("M107").Value = 1320.3333

$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H16").Value = 2561.8
$ws.Range("I16").Value = 2952.25
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 2952.25
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -2665.25
$ws.Range("N16").Value = -1574
$ws.Range("H22").Value = 822.7143
$ws.Range("I22").Value = 666.3200000000001
$ws.Range("J22").Value = 1052.7059
$ws.Range("K22").Value = 666.3200000000001
$ws.Range("L22").Value = 1052.7059
$ws.Range("M22").Value = -316.3200000000001
$ws.Range("N22").Value = -1752.7059
$ws.Range("H31").Value = 47449.76
$ws.Range("I31").Value = 3857.75
$ws.Range("J31").Value = 87688.53999999999
$ws.Range("K31").Value = 3857.75
$ws.Range("L31").Value = 87688.53999999999
$ws.Range("M31").Value = -3562.75
$ws.Range("N31").Value = -88278.53999999999
$ws.Range("H34").Value = 47449.76
$ws.Range("I34").Value = 3857.75
$ws.Range("J34").Value = 87688.53999999999
$ws.Range("K34").Value = 3857.75
$ws.Range("L34").Value = 87688.53999999999
$ws.Range("M34").Value = -3655.75
$ws.Range("N34").Value = -88092.53999999999
$ws.Range("H80").Value = 45000
$ws.Range("J80").Value = 45000
$ws.Range("L80").Value = 45000
$ws.Range("N80").Value = -47246
$ws.Range("H83").Value = 45000
$ws.Range("J83").Value = 45000
$ws.Range("L83").Value = 135000
$ws.Range("N83").Value = -146232
$ws.Range("H99").Value = 2926
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("H113").Value = 2561.8
$ws.Range("I113").Value = 2952.25
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 2952.25
$ws.Range("L113").Value = 1000
$ws.Range("M113").Value = -782.25
$ws.Range("N113").Value = -5340
$ws.Range("H126").Value = 2926
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("H134").Value = 2828
$ws.Range("I134").Value = 2333.55
$ws.Range("K134").Value = 7000.650000000001
$ws.Range("M134").Value = -4465.650000000001
$ws.Range("N99").ClearContents()
$ws.Range("N126").ClearContents()
$ws.Range("N127").ClearContents()

$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H5").Value = 9571.736999999999
$ws.Range("I5").Value = 2705.6667
$ws.Range("J5").Value = 15751.2
$ws.Range("K5").Value = 8117.000100000001
$ws.Range("L5").Value = 47253.60000000001
$ws.Range("M5").Value = -8005.000100000001
$ws.Range("N5").Value = -47477.60000000001
$ws.Range("H12").Value = 343.4
$ws.Range("I12").Value = 120.833336
$ws.Range("J12").Value = 438.7857
$ws.Range("K12").Value = 362.500008
$ws.Range("L12").Value = 1316.3571
$ws.Range("M12").Value = -189.500008
$ws.Range("N12").Value = -1662.3571
$ws.Range("H23").Value = 472.5
$ws.Range("I23").Value = 100
$ws.Range("K23").Value = 300
$ws.Range("M23").Value = -65
$ws.Range("H129").Value = 4903938
$ws.Range("I129").Value = 639.3333
$ws.Range("J129").Value = 10420149
$ws.Range("K129").Value = 1917.9999
$ws.Range("L129").Value = 31260447
$ws.Range("M129").Value = 3082.0001
$ws.Range("N129").Value = -31270447
$ws.Range("H131").Value = 6946092
$ws.Range("I131").Value = 1921.6
$ws.Range("J131").Value = 10803964
$ws.Range("K131").Value = 5764.799999999999
$ws.Range("L131").Value = 32411892
$ws.Range("M131").Value = -724.7999999999993
$ws.Range("N131").Value = -32421972
$ws.Range("H135").Value = 9571.736999999999
$ws.Range("I135").Value = 2705.6667
$ws.Range("J135").Value = 15751.2
$ws.Range("K135").Value = 24351.0003
$ws.Range("L135").Value = 141760.8
$ws.Range("M135").Value = -21816.0003
$ws.Range("N135").Value = -146830.8

$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H15").Value = 51444
$ws.Range("J15").Value = 51444
$ws.Range("L15").Value = 51444
$ws.Range("N15").Value = -52020
$ws.Range("H81").Value = 51444
$ws.Range("J81").Value = 51444
$ws.Range("L81").Value = 51444
$ws.Range("N81").Value = -53440
$ws.Range("H84").Value = 51444
$ws.Range("J84").Value = 51444
$ws.Range("L84").Value = 154332
$ws.Range("N84").Value = -164316
$ws.Range("H102").Value = 3737.3215
$ws.Range("I102").Value = 2401.5789
$ws.Range("J102").Value = 6557.222
$ws.Range("K102").Value = 2401.5789
$ws.Range("L102").Value = 6557.222
$ws.Range("M102").Value = -779.5789
$ws.Range("N102").Value = -9801.222
$ws.Range("H122").Value = 5853.4165
$ws.Range("J122").Value = 10602.667
$ws.Range("L122").Value = 31808.001
$ws.Range("N122").Value = -36708.001
$ws.Range("H132").Value = 21481.736
$ws.Range("I132").Value = 32282.605
$ws.Range("J132").Value = 3660.3
$ws.Range("K132").Value = 96847.815
$ws.Range("L132").Value = 10980.9
$ws.Range("M132").Value = -94317.815
$ws.Range("N132").Value = -16040.9

$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H40").Value = 6058.6523
$ws.Range("I40").Value = 4010.4707
$ws.Range("K40").Value = 4010.4707
$ws.Range("M40").Value = -3874.4707
$ws.Range("H61").Value = 2270.6667
$ws.Range("I61").Value = 2361.2222
$ws.Range("J61").Value = 1999
$ws.Range("K61").Value = 2361.2222
$ws.Range("L61").Value = 1999
$ws.Range("M61").Value = -2159.2222
$ws.Range("N61").Value = -2403
$ws.Range("H80").Value = 50555.445
$ws.Range("J80").Value = 50555.445
$ws.Range("L80").Value = 50555.445
$ws.Range("N80").Value = -52801.445
$ws.Range("H83").Value = 50555.445
$ws.Range("J83").Value = 50555.445
$ws.Range("L83").Value = 151666.335
$ws.Range("N83").Value = -162898.335
$ws.Range("H100").Value = 2818.1072
$ws.Range("I100").Value = 1924.0952
$ws.Range("J100").Value = 5500.143
$ws.Range("K100").Value = 1924.0952
$ws.Range("L100").Value = 5500.143
$ws.Range("M100").Value = -1383.0952
$ws.Range("N100").Value = -6582.143
$ws.Range("H110").Value = 19625
$ws.Range("J110").Value = 19625
$ws.Range("L110").Value = 19625
$ws.Range("N110").Value = -27805
$ws.Range("H113").Value = 2270.6667
$ws.Range("I113").Value = 2361.2222
$ws.Range("J113").Value = 1999
$ws.Range("K113").Value = 2361.2222
$ws.Range("L113").Value = 1999
$ws.Range("M113").Value = -191.2222000000002
$ws.Range("N113").Value = -6339
$ws.Range("H132").Value = 4740.393
$ws.Range("I132").Value = 2986.4443
$ws.Range("J132").Value = 5571.2104
$ws.Range("K132").Value = 8959.332900000001
$ws.Range("L132").Value = 16713.6312
$ws.Range("M132").Value = -6429.332900000001
$ws.Range("N132").Value = -21773.6312

$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H62").Value = 5735.5
$ws.Range("I62").Value = 2824.5
$ws.Range("K62").Value = 2824.5
$ws.Range("M62").Value = -2200.5
$ws.Range("H65").Value = 5735.5
$ws.Range("I65").Value = 2824.5
$ws.Range("K65").Value = 14122.5
$ws.Range("M65").Value = -11002.5
$ws.Range("H81").Value = 2722.739
$ws.Range("I81").Value = 1801.0555
$ws.Range("K81").Value = 3602.111
$ws.Range("M81").Value = -2541.111
$ws.Range("H84").Value = 2722.739
$ws.Range("I84").Value = 1801.0555
$ws.Range("K84").Value = 18010.555
$ws.Range("M84").Value = -12706.555
$ws.Range("H122").Value = 3257.4614
$ws.Range("I122").Value = 2386
$ws.Range("K122").Value = 7158
$ws.Range("M122").Value = -4708
$ws.Range("H126").Value = 9669.333000000001
$ws.Range("I126").Value = 8504
$ws.Range("K126").Value = 25512
$ws.Range("M126").Value = -23042
